$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.59539503035409
$ws.Range("D2").Value = 7.425240401560876
$ws.Range("E2").Value = 11.10000132855342
$ws.Range("F2").Value = 37.98404364707886
$ws.Range("G2").Value = 3.704438909715236
$ws.Range("I2").Value = 30.40166379342887
$ws.Range("K2").Value = 14.25597574865884
$ws.Range("L2").Value = 10.64563061966733
$ws.Range("M2").Value = 15.33378006908554
$ws.Range("N2").Value = 22.57931979384975
# Row 3
$ws.Range("B3").Value = 13.52813766724618
$ws.Range("D3").Value = 7.443072367879551
$ws.Range("E3").Value = 10.99274457956728
$ws.Range("F3").Value = 37.82825894697596
$ws.Range("G3").Value = 3.707789732781416
$ws.Range("I3").Value = 30.45539037862321
$ws.Range("K3").Value = 13.97989148075112
$ws.Range("L3").Value = 10.63459830886144
$ws.Range("M3").Value = 15.32699478233824
$ws.Range("N3").Value = 22.63370563443742
# Row 4
$ws.Range("B4").Value = 13.48959213279652
$ws.Range("D4").Value = 7.454438551459846
$ws.Range("E4").Value = 10.92511618158366
$ws.Range("F4").Value = 37.74205174037004
$ws.Range("G4").Value = 3.709954873915852
$ws.Range("I4").Value = 30.49280981278403
$ws.Range("K4").Value = 13.81193930320805
$ws.Range("L4").Value = 10.62984855865553
$ws.Range("M4").Value = 15.32543397972049
$ws.Range("N4").Value = 22.66900202986703
# Row 5
$ws.Range("B5").Value = 13.4745896585629
$ws.Range("D5").Value = 7.459175741947215
$ws.Range("E5").Value = 10.89711766809773
$ws.Range("F5").Value = 37.70931832416473
$ws.Range("G5").Value = 3.710864371602292
$ws.Range("I5").Value = 30.50917147733912
$ws.Range("K5").Value = 13.74399351481963
$ws.Range("L5").Value = 10.6284235334664
$ws.Range("M5").Value = 15.32545503210514
$ws.Range("N5").Value = 22.6838647528555
# Row 6
$ws.Range("B6").Value = 13.47214147208627
$ws.Range("D6").Value = 7.459968726564778
$ws.Range("E6").Value = 10.89244218312302
$ws.Range("F6").Value = 37.70402828874263
$ws.Range("G6").Value = 3.71101703777713
$ws.Range("I6").Value = 30.51195549839925
$ws.Range("K6").Value = 13.73274413972684
$ws.Range("L6").Value = 10.62821778333426
$ws.Range("M6").Value = 15.32549826544042
$ws.Range("N6").Value = 22.68636165327663
# Row 7
$ws.Range("B7").Value = 13.48938693204359
$ws.Range("D7").Value = 7.454502011620272
$ws.Range("E7").Value = 10.92474035449586
$ws.Range("F7").Value = 37.74160055585966
$ws.Range("G7").Value = 3.709967029520209
$ws.Range("I7").Value = 30.49302596750852
$ws.Range("K7").Value = 13.81102081794945
$ws.Range("L7").Value = 10.62982727132809
$ws.Range("M7").Value = 15.32543160077124
$ws.Range("N7").Value = 22.66920053295724
# Row 8
$ws.Range("B8").Value = 13.57164212319132
$ws.Range("D8").Value = 7.431302528391217
$ws.Range("E8").Value = 11.06339076647054
$ws.Range("F8").Value = 37.9283826943748
$ws.Range("G8").Value = 3.705571974534209
$ws.Range("I8").Value = 30.4192682956159
$ws.Range("K8").Value = 14.16051848767021
$ws.Range("L8").Value = 10.64140759301273
$ws.Range("M8").Value = 15.33090100689809
$ws.Range("N8").Value = 22.59767734642544
# Row 9
$ws.Range("B9").Value = 13.75413663816949
$ws.Range("D9").Value = 7.389098349963776
$ws.Range("E9").Value = 11.32095652194929
$ws.Range("F9").Value = 38.3685060934223
$ws.Range("G9").Value = 3.697803568228094
$ws.Range("I9").Value = 30.30984434118513
$ws.Range("K9").Value = 14.85384415311962
$ws.Range("L9").Value = 10.68009649033082
$ws.Range("M9").Value = 15.36220252702526
$ws.Range("N9").Value = 22.47249617848152
# Row 10
$ws.Range("B10").Value = 13.90022804874379
$ws.Range("D10").Value = 7.360066366760922
$ws.Range("E10").Value = 11.50111784312799
$ws.Range("F10").Value = 38.73520281875779
$ws.Range("G10").Value = 3.692608258016557
$ws.Range("I10").Value = 30.2509917923844
$ws.Range("K10").Value = 15.36227216609376
$ws.Range("L10").Value = 10.71813667944457
$ws.Range("M10").Value = 15.39759761387242
$ws.Range("N10").Value = 22.38967923259313
# Row 11
$ws.Range("B11").Value = 13.96908814315302
$ws.Range("D11").Value = 7.34728139649475
$ws.Range("E11").Value = 11.58103933584695
$ws.Range("F11").Value = 38.91101350451049
$ws.Range("G11").Value = 3.690354656760261
$ws.Range("I11").Value = 30.22890938198481
$ws.Range("K11").Value = 15.59219416515389
$ws.Range("L11").Value = 10.737496479756
$ws.Range("M11").Value = 15.41635435873216
$ws.Range("N11").Value = 22.35398376137485
# Row 12
$ws.Range("B12").Value = 13.99549052329903
$ws.Range("D12").Value = 7.342500231806809
$ws.Range("E12").Value = 11.61100514852373
$ws.Range("F12").Value = 38.97884276646061
$ws.Range("G12").Value = 3.689516960126642
$ws.Range("I12").Value = 30.22122263535106
$ws.Range("K12").Value = 15.67897185993824
$ws.Range("L12").Value = 10.74511965401732
$ws.Range("M12").Value = 15.42383514488128
$ws.Range("N12").Value = 22.34075081768059
# Row 13
$ws.Range("B13").Value = 13.98979008124575
$ws.Range("D13").Value = 7.343527269454962
$ws.Range("E13").Value = 11.60456485590782
$ws.Range("F13").Value = 38.96417943763393
$ws.Range("G13").Value = 3.689696676578521
$ws.Range("I13").Value = 30.22284806173334
$ws.Range("K13").Value = 15.66029699620239
$ws.Range("L13").Value = 10.74346493747155
$ws.Range("M13").Value = 15.42220727567398
$ws.Range("N13").Value = 22.34358813595667
# Row 14
$ws.Range("B14").Value = 13.97125384833868
$ws.Range("D14").Value = 7.346886842648585
$ws.Range("E14").Value = 11.58351066253938
$ws.Range("F14").Value = 38.91656897542322
$ws.Range("G14").Value = 3.69028542496173
$ws.Range("I14").Value = 30.2282634478976
$ws.Range("K14").Value = 15.5993397187764
$ws.Range("L14").Value = 10.73811781111656
$ws.Range("M14").Value = 15.41696225456256
$ws.Range("N14").Value = 22.35288938603286
# Row 15
$ws.Range("B15").Value = 13.95994182763138
$ws.Range("D15").Value = 7.348952508649726
$ws.Range("E15").Value = 11.57057528016446
$ws.Range("F15").Value = 38.88756820341596
$ws.Range("G15").Value = 3.69064809135593
$ws.Range("I15").Value = 30.23166850988602
$ws.Range("K15").Value = 15.56196135447845
$ws.Range("L15").Value = 10.73488046488036
$ws.Range("M15").Value = 15.41379863743767
$ws.Range("N15").Value = 22.35862366638973
# Row 16
$ws.Range("B16").Value = 13.89577461839552
$ws.Range("D16").Value = 7.360910345886584
$ws.Range("E16").Value = 11.49585331303168
$ws.Range("F16").Value = 38.72389085304564
$ws.Range("G16").Value = 3.692757737185738
$ws.Range("I16").Value = 30.25252938956005
$ws.Range("K16").Value = 15.34721074376915
$ws.Range("L16").Value = 10.71691253462926
$ws.Range("M16").Value = 15.39642492939471
$ws.Range("N16").Value = 22.39205179669932
# Row 17
$ws.Range("B17").Value = 13.8570121501354
$ws.Range("D17").Value = 7.368353823499087
$ws.Range("E17").Value = 11.44948823545032
$ws.Range("F17").Value = 38.62575688622158
$ws.Range("G17").Value = 3.694079986224991
$ws.Range("I17").Value = 30.26652880055044
$ws.Range("K17").Value = 15.21505135788574
$ws.Range("L17").Value = 10.70641377988112
$ws.Range("M17").Value = 15.386444201678
$ws.Range("N17").Value = 22.41306533805103
# Row 18
$ws.Range("B18").Value = 13.83494429614912
$ws.Range("D18").Value = 7.372674844824491
$ws.Range("E18").Value = 11.42262909970199
$ws.Range("F18").Value = 38.57016239038685
$ws.Range("G18").Value = 3.694850846560575
$ws.Range("I18").Value = 30.27502227774123
$ws.Range("K18").Value = 15.13891359192661
$ws.Range("L18").Value = 10.7005688414284
$ws.Range("M18").Value = 15.38095373352134
$ws.Range("N18").Value = 22.42533800353905
# Row 19
$ws.Range("B19").Value = 13.82751210447952
$ws.Range("D19").Value = 7.374144705206693
$ws.Range("E19").Value = 11.41350241375847
$ws.Range("F19").Value = 38.55148615305649
$ws.Range("G19").Value = 3.695113624931238
$ws.Range("I19").Value = 30.27797379962227
$ws.Range("K19").Value = 15.11311639468005
$ws.Range("L19").Value = 10.69862320926402
$ws.Range("M19").Value = 15.3791378355984
$ws.Range("N19").Value = 22.42952531318662
# Row 20
$ws.Range("B20").Value = 13.86111509236543
$ws.Range("D20").Value = 7.367557343878533
$ws.Range("E20").Value = 11.4544437071802
$ws.Range("F20").Value = 38.63611577270137
$ws.Range("G20").Value = 3.693938161317986
$ws.Range("I20").Value = 30.26499285019116
$ws.Range("K20").Value = 15.22913338498192
$ws.Range("L20").Value = 10.70751137067126
$ws.Range("M20").Value = 15.38748080096736
$ws.Range("N20").Value = 22.41080913893643
# Row 21
$ws.Range("B21").Value = 13.97668967989554
$ws.Range("D21").Value = 7.345898423480174
$ws.Range("E21").Value = 11.58970294461419
$ws.Range("F21").Value = 38.93051963484341
$ws.Range("G21").Value = 3.69011207010819
$ws.Range("I21").Value = 30.22665448085868
$ws.Range("K21").Value = 15.61725289970532
$ws.Range("L21").Value = 10.73968049527939
$ws.Range("M21").Value = 15.41849261723735
$ws.Range("N21").Value = 22.35014967401809
# Row 22
$ws.Range("B22").Value = 14.05411865217032
$ws.Range("D22").Value = 7.332093911040887
$ws.Range("E22").Value = 11.67635878509456
$ws.Range("F22").Value = 39.13021592826168
$ws.Range("G22").Value = 3.687702928389324
$ws.Range("I22").Value = 30.20553514885865
$ws.Range("K22").Value = 15.86918878993001
$ws.Range("L22").Value = 10.76240525895772
$ws.Range("M22").Value = 15.44096199145269
$ws.Range("N22").Value = 22.31216135113647
# Row 23
$ws.Range("B23").Value = 14.0126264021503
$ws.Range("D23").Value = 7.339429679234659
$ws.Range("E23").Value = 11.63027047319752
$ws.Range("F23").Value = 39.02298152578198
$ws.Range("G23").Value = 3.688980396988975
$ws.Range("I23").Value = 30.21644640798709
$ws.Range("K23").Value = 15.73491279703276
$ws.Range("L23").Value = 10.75012227358749
$ws.Range("M23").Value = 15.42876958295288
$ws.Range("N23").Value = 22.33228500029265
# Row 24
$ws.Range("B24").Value = 13.85925947419189
$ws.Range("D24").Value = 7.367917302479338
$ws.Range("E24").Value = 11.45220397105108
$ws.Range("F24").Value = 38.63142995027771
$ws.Range("G24").Value = 3.694002247078662
$ws.Range("I24").Value = 30.26568586717849
$ws.Range("K24").Value = 15.22276738423666
$ws.Range("L24").Value = 10.70701455487193
$ws.Range("M24").Value = 15.38701138294605
$ws.Range("N24").Value = 22.41182856937432
# Row 25
$ws.Range("B25").Value = 13.70258865388407
$ws.Range("D25").Value = 7.400166614865701
$ws.Range("E25").Value = 11.25284737729197
$ws.Range("F25").Value = 38.2416967833429
$ws.Range("G25").Value = 3.699814743390977
$ws.Range("I25").Value = 30.33566841090906
$ws.Range("K25").Value = 14.6660403105973
$ws.Range("L25").Value = 10.6679310918446
$ws.Range("M25").Value = 15.35154561508248
$ws.Range("N25").Value = 22.5047509416405
